# PrefixeSocoda : ajout de nouveaux fournisseurs dans les prefixes / correction
# nom mediator / ajout d'une colonne "temps" pour MERSEN.
#
# Resume des modifications :
#   - D33 (ligne MERSEN) : ajout de la duree "19m40"
#   - Nouvelle ligne fournisseur "TE CONNECTIVITY" / "ENTRELEC" / "ENT" avec le
#     commentaire "pas de socoda" (inseree a sa place alphabetique, juste avant
#     "TECTHRONIC INDUSTRIES FRANCE", donc toutes les lignes suivantes sont
#     decalees d'un cran vers le bas)
#   - Suppression du filtre automatique (AutoFilter)
#   - Re-application du tri sur la plage mise a jour

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Ajoute la duree manquante sur la ligne MERSEN existante (ligne 33)
$ws.Range("D33").Value2 = "19m40"

# 2) Insere une nouvelle ligne 44 (ce qui decale vers le bas les lignes
#    44 a 51 existantes, qui deviennent 45 a 52) et saisit les informations du
#    nouveau fournisseur TE CONNECTIVITY / ENTRELEC
$ws.Rows("44:44").Insert()
$ws.Range("B44").Value2 = "ENTRELEC"
$ws.Range("A44").Value2 = "TE CONNECTIVITY"
$ws.Range("F44").Value2 = "pas de socoda"
$ws.Range("C44").Value2 = "ENT"

# 3) Retire le filtre automatique de la feuille
$ws.AutoFilterMode = $false

# 4) Met a jour l'etat de tri (sortState) pour refleter la nouvelle plage de
#    donnees, triee par la colonne A (FABRICANT)
$sortObj = $ws.Sort
$sortObj.SetRange($ws.Range("A1:F52"))
$sortObj.SortFields.Clear()
$sortObj.SortFields.Add($ws.Range("A2:A52"))
$sortObj.Header = 1
$sortObj.Apply()

# 5) Replace la vue / selection telle qu'elle etait au moment de
#    l'enregistrement (cosmetique)
$win = $excel.ActiveWindow
$win.ScrollRow = 32
$win.ScrollColumn = 1
$ws.Range("C55:D55").Select() | Out-Null
